$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Login with valid username and password"
$ws.Range("B8").Value = "Passed"

$ws.Range("A9").Value = "Create a Citizenship with parameter"
$ws.Range("B9").Value = "Failed"
